# Applies the "feat: add 2022-Q1 data" change:
#  - Inserts a new worksheet "2022-Q1" (with per-fund holding detail) right
#    before the "总计" (summary) sheet.
#  - Updates the "总计" sheet by adding a new top data row for 2022-Q1 and
#    shifting the previous rows down.
#
# To obtain the exact sheetId/rId numbering shown in the target diff
# (2022-Q1 -> sheetId 6 / rId6, 总计 -> sheetId 7 / rId7) the existing
# "总计" sheet is removed and both sheets are (re)created in order.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value as literal text (keeps numeric-looking strings such
# as "513700" or "3.24" stored as text, matching the workbook convention).
# ---------------------------------------------------------------------------
function Set-TextValue($ws, $row, $col, [string]$val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
}

function Set-NumberValue($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

function Format-HeaderCell($ws, $row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# Locate and remove the existing "总计" sheet (last sheet in the workbook);
# we recreate it after the new "2022-Q1" sheet so sheetId/rId numbering
# lines up with the target workbook.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalSheet.Delete()

$anchor = $wb.Worksheets.Item($wb.Worksheets.Count)

# ---------------------------------------------------------------------------
# New "2022-Q1" sheet: per-fund holding detail.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $anchor)
$q1.Name = "2022-Q1"

$q1.Outline.SummaryRow = 1
$q1.Outline.SummaryColumn = 1
$q1.PageSetup.LeftMargin = 54
$q1.PageSetup.RightMargin = 54
$q1.PageSetup.TopMargin = 72
$q1.PageSetup.BottomMargin = 72
$q1.PageSetup.HeaderMargin = 36
$q1.PageSetup.FooterMargin = 36

# Header row
Set-TextValue $q1 1 2 "基金代码"
Set-TextValue $q1 1 3 "基金名称"
Set-TextValue $q1 1 4 "基金规模"
Set-TextValue $q1 1 5 "股票总仓位"
Set-TextValue $q1 1 6 "仓位占比"
Set-TextValue $q1 1 7 "持有市值(亿元)"
Set-TextValue $q1 1 8 "仓位排名"
foreach ($col in 2..8) { Format-HeaderCell $q1 1 $col }

$q1Rows = @(
    @("513700", "鹏华中证港股通医药卫生综合交易型开放式指数证券投资基金", "3.24", "93.11", "8.26", "0.2676", 2),
    @("513980", "景顺长城中证港股通科技交易型开放式指数证券投资基金",     "5.03", "97.36", "4.06", "0.2042", 7),
    @("517120", "华泰柏瑞中证沪港深创新药产业交易型开放式指数证券投资基金", "4.15", "95.09", "4.75", "0.1971", 6),
    @("513860", "海富通中证港股通科技交易型开放式指数证券投资基金",       "3.87", "95.13", "4.18", "0.1618", 7),
    @("862001", "光大阳光香港精选混合型集合资产管理计划（QDII）A 人民币", "3.91", "89.45", "4.00", "0.1564", 9),
    @("862011", "光大阳光香港精选混合型集合资产管理计划（QDII）A 美元",   "3.91", "89.45", "4.00", "0.1564", 9),
    @("862012", "光大阳光香港精选混合型集合资产管理计划（QDII）C 人民币", "3.91", "89.45", "4.00", "0.1564", 9),
    @("159748", "富国中证沪港深创新药产业ETF",                         "1.94", "99.28", "4.78", "0.0927", 6),
    @("159747", "南方中证香港科技交易型开放式指数证券投资基金(QDII)",     "3.07", "100.03", "2.97", "0.0912", 9),
    @("517110", "国泰中证沪港深创新药产业ETF",                         "1.66", "94.10", "3.91", "0.0649", 6),
    @("159751", "鹏华中证港股通科技ETF",                                "0.90", "91.30", "3.88", "0.0349", 7),
    @("513020", "国泰中证港股通科技ETF",                                "0.76", "91.35", "3.60", "0.0274", 7),
    @("006786", "泰康中证港股通大消费主题指数A",                        "0.85", "80.77", "2.50", "0.0212", 10),
    @("860008", "光大阳光生活 18 个月持有期混合型集合资产管理计划A",      "0.51", "88.26", "3.13", "0.0160", 10),
    @("006787", "泰康中证港股通大消费主题指数C",                        "0.41", "80.77", "2.50", "0.0102", 10),
    @("860060", "光大阳光生活 18 个月持有期混合型集合资产管理计划B",      "0.14", "88.26", "3.13", "0.0044", 10),
    @("860061", "光大阳光生活 18 个月持有期混合型集合资产管理计划C",      "0.03", "88.26", "3.13", "0.0009", 10)
)

$r = 2
foreach ($row in $q1Rows) {
    $idx = $r - 2
    Set-NumberValue $q1 $r 1 $idx
    Set-TextValue   $q1 $r 2 $row[0]
    Set-TextValue   $q1 $r 3 $row[1]
    Set-TextValue   $q1 $r 4 $row[2]
    Set-TextValue   $q1 $r 5 $row[3]
    Set-TextValue   $q1 $r 6 $row[4]
    Set-TextValue   $q1 $r 7 $row[5]
    Set-NumberValue $q1 $r 8 $row[6]
    $q1.Cells.Item($r, 1).Font.Bold = $true
    $q1.Cells.Item($r, 1).Borders.LineStyle = 1
    $q1.Cells.Item($r, 1).HorizontalAlignment = -4108
    $q1.Cells.Item($r, 1).VerticalAlignment = -4160
    $r++
}

# ---------------------------------------------------------------------------
# Recreate the "总计" (summary) sheet with the new 2022-Q1 row prepended.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.Outline.SummaryRow = 1
$total.Outline.SummaryColumn = 1
$total.PageSetup.LeftMargin = 54
$total.PageSetup.RightMargin = 54
$total.PageSetup.TopMargin = 72
$total.PageSetup.BottomMargin = 72
$total.PageSetup.HeaderMargin = 36
$total.PageSetup.FooterMargin = 36

Set-TextValue $total 1 2 "日期"
Set-TextValue $total 1 3 "持有数量(只)"
Set-TextValue $total 1 4 "持有市值(亿元)"
foreach ($col in 2..4) { Format-HeaderCell $total 1 $col }

$totalRows = @(
    @("2022-Q1", 17, 1.66),
    @("2021-Q4", 8, 1.4),
    @("2021-Q3", 12, 2.17),
    @("2021-Q2", 5, 1.03),
    @("2021-Q1", 2, 0.06),
    @("2020-Q4", 2, 0.03)
)

$r = 2
foreach ($row in $totalRows) {
    $idx = $r - 2
    Set-NumberValue $total $r 1 $idx
    Set-TextValue   $total $r 2 $row[0]
    Set-NumberValue $total $r 3 $row[1]
    Set-NumberValue $total $r 4 $row[2]
    $total.Cells.Item($r, 1).Font.Bold = $true
    $total.Cells.Item($r, 1).Borders.LineStyle = 1
    $total.Cells.Item($r, 1).HorizontalAlignment = -4108
    $total.Cells.Item($r, 1).VerticalAlignment = -4160
    $r++
}
